$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 6339.902
$ws.Range("I40").Value = 5043.4546
$ws.Range("J40").Value = 8716.723
$ws.Range("K40").Value = 5043.4546
$ws.Range("L40").Value = 8716.723
$ws.Range("M40").Value = -4868.4546
$ws.Range("N40").Value = -9066.723
# Row 75
$ws.Range("H75").Value = 60314
$ws.Range("J75").Value = 60314
$ws.Range("L75").Value = 60314
$ws.Range("N75").Value = -62186
# Row 78
$ws.Range("H78").Value = 60314
$ws.Range("J78").Value = 60314
$ws.Range("L78").Value = 180942
$ws.Range("N78").Value = -190302
# Row 107
$ws.Range("H107").Value = 382.77777
$ws.Range("I107").Value = 182.83333
$ws.Range("J107").Value = 782.6667
$ws.Range("K107").Value = 182.83333
$ws.Range("L107").Value = 782.6667
$ws.Range("M107").Value = 1737.16667
$ws.Range("N107").Value = -4622.6667
# Row 112
$ws.Range("H112").Value = 2137.5
$ws.Range("I112").Value = 2200
$ws.Range("J112").Value = 2121.875
$ws.Range("K112").Value = 6600
$ws.Range("L112").Value = 6365.625
$ws.Range("M112").Value = -5492
$ws.Range("N112").Value = -8581.625
# Row 125
$ws.Range("H125").Value = 1675.25
$ws.Range("I125").Value = 1933.6666
$ws.Range("K125").Value = 17402.9994
$ws.Range("M125").Value = -14942.9994
# Row 132
$ws.Range("H132").Value = 21429.3
$ws.Range("I132").Value = 21429.3
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 64287.89999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -61757.89999999999
$ws.Range("N132").ClearContents()
# Row 138
$ws.Range("H138").Value = 5739.4
$ws.Range("I138").Value = 1232.3334
$ws.Range("J138").Value = 12500
$ws.Range("K138").Value = 3697.0002
$ws.Range("L138").Value = 37500
$ws.Range("M138").Value = 1442.9998
$ws.Range("N138").Value = -47780

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 79.53846
$ws.Range("I5").Value = 73.40000000000001
$ws.Range("K5").Value = 73.40000000000001
$ws.Range("M5").Value = 38.59999999999999
# Row 32
$ws.Range("H32").Value = 820.9091
$ws.Range("I32").Value = 529.29266
$ws.Range("J32").Value = 4806.3335
$ws.Range("K32").Value = 529.29266
$ws.Range("L32").Value = 4806.3335
$ws.Range("M32").Value = -242.29266
$ws.Range("N32").Value = -5380.3335
# Row 122
$ws.Range("H122").Value = 3331.6667
$ws.Range("I122").Value = 3598
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 10794
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -8344
$ws.Range("N122").Value = -10900
# Row 132
$ws.Range("H132").Value = 4576.857
$ws.Range("I132").Value = 4589.6665
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 13768.9995
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -11238.9995
$ws.Range("N132").Value = -18560

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 79.53846
$ws.Range("I4").Value = 73.40000000000001
$ws.Range("K4").Value = 73.40000000000001
$ws.Range("M4").Value = 41.59999999999999
# Row 80
$ws.Range("H80").Value = 1425.7142
$ws.Range("I80").Value = 1071.6666
$ws.Range("J80").Value = 1691.25
$ws.Range("K80").Value = 1071.6666
$ws.Range("L80").Value = 1691.25
$ws.Range("M80").Value = -73.66660000000002
$ws.Range("N80").Value = -3687.25
# Row 83
$ws.Range("H83").Value = 1425.7142
$ws.Range("I83").Value = 1071.6666
$ws.Range("J83").Value = 1691.25
$ws.Range("K83").Value = 5358.333000000001
$ws.Range("L83").Value = 8456.25
$ws.Range("M83").Value = -366.3330000000005
$ws.Range("N83").Value = -18440.25
# Row 86
$ws.Range("H86").Value = 3721.1875
$ws.Range("I86").Value = 2310.889
$ws.Range("J86").Value = 5534.4287
$ws.Range("K86").Value = 2310.889
$ws.Range("L86").Value = 5534.4287
$ws.Range("M86").Value = -1187.889
$ws.Range("N86").Value = -7780.4287
# Row 89
$ws.Range("H89").Value = 3721.1875
$ws.Range("I89").Value = 2310.889
$ws.Range("J89").Value = 5534.4287
$ws.Range("K89").Value = 11554.445
$ws.Range("L89").Value = 27672.1435
$ws.Range("M89").Value = -5938.445
$ws.Range("N89").Value = -38904.14350000001
# Row 99
$ws.Range("H99").Value = 2127.25
$ws.Range("I99").Value = 2127.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2127.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -629.25
$ws.Range("N99").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6043.9644
$ws.Range("I31").Value = 2506.7693
$ws.Range("J31").Value = 9109.532999999999
$ws.Range("K31").Value = 2506.7693
$ws.Range("L31").Value = 9109.532999999999
$ws.Range("M31").Value = -2211.7693
$ws.Range("N31").Value = -9699.532999999999
# Row 34
$ws.Range("H34").Value = 6043.9644
$ws.Range("I34").Value = 2506.7693
$ws.Range("J34").Value = 9109.532999999999
$ws.Range("K34").Value = 2506.7693
$ws.Range("L34").Value = 9109.532999999999
$ws.Range("M34").Value = -2304.7693
$ws.Range("N34").Value = -9513.532999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 207.5
$ws.Range("J22").Value = 207.5
$ws.Range("L22").Value = 622.5
$ws.Range("N22").Value = -960.5
# Row 23
$ws.Range("H23").Value = 277.55554
$ws.Range("I23").Value = 329.8
$ws.Range("J23").Value = 212.25
$ws.Range("K23").Value = 989.4000000000001
$ws.Range("L23").Value = 636.75
$ws.Range("M23").Value = -754.4000000000001
$ws.Range("N23").Value = -1106.75
# Row 27
$ws.Range("H27").Value = 207.5
$ws.Range("J27").Value = 207.5
$ws.Range("L27").Value = 622.5
$ws.Range("N27").Value = -826.5
# Row 68
$ws.Range("H68").Value = 3333
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3333
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 9999
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11621
# Row 71
$ws.Range("H71").Value = 3333
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3333
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 29997
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -38109
# Row 112
$ws.Range("H112").Value = 625.5
$ws.Range("I112").Value = 625.5
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 1876.5
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -768.5
$ws.Range("N112").ClearContents()
# Row 117
$ws.Range("H117").Value = 139.875
$ws.Range("J117").Value = 149.85715
$ws.Range("L117").Value = 449.57145
$ws.Range("N117").Value = -7333.57145

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 3909.2856
$ws.Range("I70").Value = 3978.6
$ws.Range("J70").Value = 3736
$ws.Range("K70").Value = 3978.6
$ws.Range("L70").Value = 3736
$ws.Range("M70").Value = -3708.6
$ws.Range("N70").Value = -4276
# Row 73
$ws.Range("H73").Value = 3909.2856
$ws.Range("I73").Value = 3978.6
$ws.Range("J73").Value = 3736
$ws.Range("K73").Value = 3978.6
$ws.Range("L73").Value = 3736
$ws.Range("M73").Value = -3042.6
$ws.Range("N73").Value = -5608
# Row 122
$ws.Range("H122").Value = 3066
$ws.Range("I122").Value = 2666.5
$ws.Range("J122").Value = 3332.3333
$ws.Range("K122").Value = 7999.5
$ws.Range("L122").Value = 9996.999899999999
$ws.Range("M122").Value = -5549.5
$ws.Range("N122").Value = -14896.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 130
$ws.Range("H130").Value = 78000
$ws.Range("J130").Value = 78000
$ws.Range("L130").Value = 78000
$ws.Range("N130").Value = -88040
# Row 132
$ws.Range("H132").Value = 3869.3333
$ws.Range("I132").Value = 3869.3333
$ws.Range("K132").Value = 11607.9999
$ws.Range("M132").Value = -9077.999899999999
# Row 136
$ws.Range("H136").Value = 7499.5
$ws.Range("I136").Value = 7499
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 22497
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -19947
$ws.Range("N136").Value = -27600
